# New weekly data point: insert two rows (Primera / Segunda) at the top of
# the Perejil data block (row 92), pushing the existing history down by two
# rows. The new rows duplicate the price/quality pattern of the rows
# immediately below them (now at 94:95, after the insert) but are dated
# with the newest report date (serial 45238 = 2023-11-08).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 92 - everything currently at
# row 92 and below shifts down to row 94 and below.
$ws.Rows("92:93").Insert()

# Seed the two new rows with the same values as the pair that is now just
# below them (rows 94:95 - the old rows 92:93), then overwrite the date.
$ws.Range("A94:R95").Copy()
$ws.Range("A92").PasteSpecial()

$ws.Range("D92").Value = 45238
$ws.Range("D93").Value = 45238
